$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45014
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("S2").Value = 1000

# Row 3
$ws.Range("D3").Value = 45001
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17500
$ws.Range("S3").Value = 972

# Row 4
$ws.Range("D4").Value = 44999
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("S4").Value = 972

# Row 5 - update with new data (Provincia de Los Andes entry)
$ws.Range("D5").Value = 45020
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "`$/caja 16 kilos"
$ws.Range("R5").Value = "Provincia de Los Andes"
$ws.Range("S5").Value = 938
$ws.Range("T5").Value = 16

# Row 6 - new row, copy of old row3/4 style data (Región Metropolitana)
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 45002
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107011
$ws.Range("J6").Value = "Tuna"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = "`$/caja 18 kilos"
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 18

# Match the date cell style used in D2:D5 (numFmt index style 2)
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
